$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "София, 2020" -> "София, 2021"
#    Only the last run (containing just "0") changes to "1"; the
#    paragraph's other two runs ("София, 20" and "2") are untouched.
# ------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("София, 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $lastDigit = $d.Range($r1.End - 1, $r1.End)
    $lastDigit.Text = "1"
}

# ------------------------------------------------------------------
# 2) Hyperlink text split across two runs ("...ondemand.com" + "/")
#    gets combined into a single run.
# ------------------------------------------------------------------
$r2 = $d.Content
$hyperlinkText = "https://vvps-project-test-report.cfapps.sap.hana.ondemand.com/"
$r2.Find.Execute($hyperlinkText, $true, $false, $false, $false, $false, $true, 1, $false, $hyperlinkText, 2) | Out-Null

# ------------------------------------------------------------------
# 3) Heading "Доклад за извършените " + "функционални" + " тестове"
#    (three runs) gets combined into a single run.
# ------------------------------------------------------------------
$r3 = $d.Content
$headingText = "Доклад за извършените функционални тестове"
$r3.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, $headingText, 2) | Out-Null

# ------------------------------------------------------------------
# 4) Dispersion paragraph: two separate places where adjacent runs
#    with identical formatting get combined into a single run.
#    The second merge ("компонент за потребител" + ", а третият ...")
#    must happen before the first ("...дисперсията " + "на броя
#    редакции на ") to avoid unwanted cascading merges elsewhere in
#    the paragraph.
# ------------------------------------------------------------------
$r4a = $d.Content
$tailText = "компонент за потребител, а третият тест валидира стандартното отклонение."
$r4a.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, $tailText, 2) | Out-Null

$r4b = $d.Content
$midText = "компонент за потребител. Вторият тест валидира правилното изчисление на дисперсията на броя редакции на "
$r4b.Find.Execute($midText, $true, $false, $false, $false, $false, $true, 1, $false, $midText, 2) | Out-Null
